$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-24 Wednesday", "2025-12-25 Thursday"),
    @("49×58=", "86×23="),
    @("46×70=", "20×29="),
    @("77×70=", "96×59="),
    @("65×16=", "89×22="),
    @("72×62=", "85×50="),
    @("61×51=", "87×90="),
    @("95×99=", "17×30="),
    @("73×71=", "99×58="),
    @("12×71=", "15×80="),
    @("55×22=", "98×47="),
    @("24×98=", "93×99="),
    @("21×42=", "97×32="),
    @("66×69=", "15×95="),
    @("40×68=", "26×28="),
    @("43×42=", "53×58="),
    @("50×19=", "14×37="),
    @("74×62=", "59×86="),
    @("55×85=", "26×70="),
    @("93×66=", "96×22="),
    @("30×81=", "64×98="),
    @("27×69=", "70×47="),
    @("29×96=", "99×41="),
    @("22×14=", "66×90="),
    @("62×12=", "32×59="),
    @("41×39=", "70×79=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
